$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B header
$ws.Range("B1").Value = "y"

# B2: formula referencing A2 (text "NA")
$ws.Range("B2").Formula = "=A2"

# B3:B4 share a formula referencing the same row in column A
$ws.Range("B3:B4").Formula = "=A3"

# Move the active selection to B2
$ws.Range("B2").Select() | Out-Null
